$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.799319333333333
$ws.Range("H2").Value = 11.397958
$ws.Range("I2").Value = 0.04516380397110807
$ws.Range("J2").Value = 0.04516380397110807
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 605.9292666708552
$ws.Range("R2").Value = 5453.363400037697
$ws.Range("S2").Value = 0.01347365027644561
$ws.Range("T2").Value = 0.01347365027644561
$ws.Range("G3").Value = 3.799319333333333
$ws.Range("H3").Value = 11.397958
$ws.Range("I3").Value = 0.04516380397110807
$ws.Range("J3").Value = 0.04516380397110807
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 655.604833783042
$ws.Range("R3").Value = 5900.443504047378
$ws.Range("S3").Value = 0.01457825316554369
$ws.Range("T3").Value = 0.01457825316554369
$ws.Range("G4").Value = 3.799319333333333
$ws.Range("H4").Value = 11.397958
$ws.Range("I4").Value = 0.04516380397110807
$ws.Range("J4").Value = 0.04516380397110807
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 282.6226483016762
$ws.Range("R4").Value = 2543.603834715086
$ws.Range("S4").Value = 0.006284493806251775
$ws.Range("T4").Value = 0.006284493806251776
$ws.Range("G5").Value = 3.799319333333333
$ws.Range("H5").Value = 11.397958
$ws.Range("I5").Value = 0.04516380397110807
$ws.Range("J5").Value = 0.04516380397110807
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 221.9357735858438
$ws.Range("R5").Value = 1997.421962272594
$ws.Range("S5").Value = 0.004935039717684435
$ws.Range("T5").Value = 0.004935039717684436
$ws.Range("G6").Value = 3.799319333333333
$ws.Range("H6").Value = 11.397958
$ws.Range("I6").Value = 0.04516380397110807
$ws.Range("J6").Value = 0.04516380397110807
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 264.9881468756427
$ws.Range("R6").Value = 2384.893321880784
$ws.Range("S6").Value = 0.005892367005182561
$ws.Range("T6").Value = 0.005892367005182563
$ws.Range("I7").Value = 0.7704071060399821
$ws.Range("J7").Value = 0.770407106039982
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 10335.98084650816
$ws.Range("R7").Value = 93023.82761857341
$ws.Range("S7").Value = 0.2298344028751791
$ws.Range("T7").Value = 0.229834402875179
$ws.Range("I8").Value = 0.7704071060399821
$ws.Range("J8").Value = 0.770407106039982
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.2486767908117189
$ws.Range("T8").Value = 0.2486767908117189
$ws.Range("I9").Value = 0.7704071060399821
$ws.Range("J9").Value = 0.770407106039982
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 4820.995519304306
$ws.Range("R9").Value = 43388.95967373875
$ws.Range("S9").Value = 0.1072013041527209
$ws.Range("T9").Value = 0.1072013041527209
$ws.Range("I10").Value = 0.7704071060399821
$ws.Range("J10").Value = 0.770407106039982
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 3785.794862726654
$ws.Range("R10").Value = 34072.15376453989
$ws.Range("S10").Value = 0.08418222852808906
$ws.Range("T10").Value = 0.08418222852808906
$ws.Range("I11").Value = 0.7704071060399821
$ws.Range("J11").Value = 0.770407106039982
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 4520.185046856515
$ws.Range("R11").Value = 40681.66542170863
$ws.Range("S11").Value = 0.1005123796722741
$ws.Range("T11").Value = 0.1005123796722741
$ws.Range("G12").Value = 2.946166666666667
$ws.Range("H12").Value = 8.8385
$ws.Range("I12").Value = 0.03502208741238024
$ws.Range("J12").Value = 0.03502208741238024
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 469.8653761902223
$ws.Range("R12").Value = 4228.788385712
$ws.Range("S12").Value = 0.01044808710194971
$ws.Range("T12").Value = 0.01044808710194971
$ws.Range("G13").Value = 2.946166666666667
$ws.Range("H13").Value = 8.8385
$ws.Range("I13").Value = 0.03502208741238024
$ws.Range("J13").Value = 0.03502208741238024
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 508.3860919114999
$ws.Range("R13").Value = 4575.474827203499
$ws.Range("S13").Value = 0.01130464690286259
$ws.Range("T13").Value = 0.01130464690286259
$ws.Range("G14").Value = 2.946166666666667
$ws.Range("H14").Value = 8.8385
$ws.Range("I14").Value = 0.03502208741238024
$ws.Range("J14").Value = 0.03502208741238024
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 219.1585788449444
$ws.Range("R14").Value = 1972.4272096045
$ws.Range("S14").Value = 0.004873285066198376
$ws.Range("T14").Value = 0.004873285066198376
$ws.Range("G15").Value = 2.946166666666667
$ws.Range("H15").Value = 8.8385
$ws.Range("I15").Value = 0.03502208741238024
$ws.Range("J15").Value = 0.03502208741238024
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 172.0991895950555
$ws.Range("R15").Value = 1548.8927063555
$ws.Range("S15").Value = 0.003826856402239232
$ws.Range("T15").Value = 0.003826856402239233
$ws.Range("G16").Value = 2.946166666666667
$ws.Range("H16").Value = 8.8385
$ws.Range("I16").Value = 0.03502208741238024
$ws.Range("J16").Value = 0.03502208741238024
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 205.4839767053333
$ws.Range("R16").Value = 1849.355790348
$ws.Range("S16").Value = 0.00456921193913033
$ws.Range("T16").Value = 0.004569211939130332
$ws.Range("G17").Value = 10.035916
$ws.Range("H17").Value = 30.107748
$ws.Range("I17").Value = 0.1193003543865946
$ws.Range("J17").Value = 0.1193003543865946
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 1600.564387651797
$ws.Range("R17").Value = 14405.07948886618
$ws.Range("S17").Value = 0.03559069678650813
$ws.Range("T17").Value = 0.03559069678650813
$ws.Range("G18").Value = 10.035916
$ws.Range("H18").Value = 30.107748
$ws.Range("I18").Value = 0.1193003543865946
$ws.Range("J18").Value = 0.1193003543865946
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 1731.782580978252
$ws.Range("R18").Value = 15586.04322880427
$ws.Range("S18").Value = 0.03850850938285539
$ws.Range("T18").Value = 0.03850850938285539
$ws.Range("G19").Value = 10.035916
$ws.Range("H19").Value = 30.107748
$ws.Range("I19").Value = 0.1193003543865946
$ws.Range("J19").Value = 0.1193003543865946
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 746.5487655033907
$ws.Range("R19").Value = 6718.938889530516
$ws.Range("S19").Value = 0.01660051351533224
$ws.Range("T19").Value = 0.01660051351533224
$ws.Range("G20").Value = 10.035916
$ws.Range("H20").Value = 30.107748
$ws.Range("I20").Value = 0.1193003543865946
$ws.Range("J20").Value = 0.1193003543865946
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 586.2441626217293
$ws.Range("R20").Value = 5276.197463595564
$ws.Range("S20").Value = 0.01303592557456644
$ws.Range("T20").Value = 0.01303592557456644
$ws.Range("G21").Value = 10.035916
$ws.Range("H21").Value = 30.107748
$ws.Range("I21").Value = 0.1193003543865946
$ws.Range("J21").Value = 0.1193003543865946
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 699.967165093856
$ws.Range("R21").Value = 6299.704485844704
$ws.Range("S21").Value = 0.01556470912733239
$ws.Range("T21").Value = 0.01556470912733239
$ws.Range("G22").Value = 2.532664666666667
$ws.Range("H22").Value = 7.597994
$ws.Range("I22").Value = 0.03010664818993502
$ws.Range("J22").Value = 0.03010664818993501
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 403.9185731856143
$ws.Range("R22").Value = 3635.267158670528
$ws.Range("S22").Value = 0.00898167145014327
$ws.Range("T22").Value = 0.008981671450143268
$ws.Range("G23").Value = 2.532664666666667
$ws.Range("H23").Value = 7.597994
$ws.Range("I23").Value = 0.03010664818993502
$ws.Range("J23").Value = 0.03010664818993501
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 437.032808285006
$ws.Range("R23").Value = 3933.295274565054
$ws.Range("S23").Value = 0.009718010900047355
$ws.Range("T23").Value = 0.009718010900047355
$ws.Range("G24").Value = 2.532664666666667
$ws.Range("H24").Value = 7.597994
$ws.Range("I24").Value = 0.03010664818993502
$ws.Range("J24").Value = 0.03010664818993501
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 188.3991137763664
$ws.Range("R24").Value = 1695.592023987298
$ws.Range("S24").Value = 0.004189307087544818
$ws.Range("T24").Value = 0.004189307087544818
$ws.Range("G25").Value = 2.532664666666667
$ws.Range("H25").Value = 7.597994
$ws.Range("I25").Value = 0.03010664818993502
$ws.Range("J25").Value = 0.03010664818993501
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 147.9446297389936
$ws.Range("R25").Value = 1331.501667650942
$ws.Range("S25").Value = 0.003289747353405587
$ws.Range("T25").Value = 0.003289747353405587
$ws.Range("G26").Value = 2.532664666666667
$ws.Range("H26").Value = 7.597994
$ws.Range("I26").Value = 0.03010664818993502
$ws.Range("J26").Value = 0.03010664818993501
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 176.6437768969014
$ws.Range("R26").Value = 1589.793992072112
$ws.Range("S26").Value = 0.003927911398793983
$ws.Range("T26").Value = 0.003927911398793983
